$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: Cons Credit - Revolving (REVOLSL) ---
# Reuse the existing "updated date" style (same one already used by N29/N30/N47-51,
# fill + yyyy-mm-dd number format) instead of inventing a new style/fill.
$ws.Range("N29").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = 45901
$ws.Range("F24").Value = 0.001263881152865798
$ws.Range("G24").Value = -0.004631607864557297
$ws.Range("H24").Value = 0.008466932005238847
$ws.Range("I24").Value = [double]"7.468713370495372e-05"
$ws.Range("J24").Value = -0.0009803765564608824

# --- Row 25: Cons Credit - NonRevolving (NONREVSL) ---
$ws.Range("N29").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 45901
$ws.Range("F25").Value = 0.003044704234780982
$ws.Range("G25").Value = 0.002454000808095547
$ws.Range("H25").Value = 0.001792092542795443
$ws.Range("I25").Value = -0.001282602659991805
$ws.Range("J25").Value = 0.002439153249045978

# --- Row 29: 5yr, 5yr Forward (T5YIFR) ---
$ws.Range("N29").Value = 45968
$ws.Range("Q29").Value = 2.2
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = 2.21
$ws.Range("T29").Value = 2.2
$ws.Range("U29").Value = 2.21

# --- Row 30: 10yr TIPS (T10YIE) ---
$ws.Range("N30").Value = 45968
$ws.Range("Q30").Value = 2.28
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.3
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.31

# --- Row 47: FFR (DFF) ---
$ws.Range("N47").Value = 45967
$ws.Range("T47").Value = 3.87

# --- Row 48: 2y UST (DGS2) ---
$ws.Range("N48").Value = 45967
$ws.Range("Q48").Value = 3.57
$ws.Range("R48").Value = 3.63
$ws.Range("S48").Value = 3.58
$ws.Range("T48").Value = 3.6

# --- Row 49: 5y UST (DGS5) ---
$ws.Range("N49").Value = 45967
$ws.Range("Q49").Value = 3.69
$ws.Range("R49").Value = 3.76
$ws.Range("S49").Value = 3.69
$ws.Range("T49").Value = 3.72

# --- Row 50: 10y UST (DGS10) ---
$ws.Range("N50").Value = 45967
$ws.Range("Q50").Value = 4.11
$ws.Range("R50").Value = 4.17
$ws.Range("S50").Value = 4.1
$ws.Range("T50").Value = 4.13

# --- Row 52: BAA (DBAA) ---
$ws.Range("N52").Value = 45967
$ws.Range("Q52").Value = 5.83
$ws.Range("R52").Value = 5.87
$ws.Range("S52").Value = 5.82
$ws.Range("T52").Value = 5.84
